# Generate Report for Handoff
# Adds two new files (a41100b7-... and be004eef-...) to the localization
# status report: one new row on the "Overview" sheet and one new row on
# each of the "zh-cn" and "de-de" sheets, keeping each sheet's Excel table
# in sync with the newly added rows.

$wb = $excel.ActiveWorkbook

$dateHandoff = "2016-08-25 12:43:26"
$dateZhCn    = "2016-08-25 12:43:22"

$file1Name  = "a41100b7-a95c-4c82-9f5a-4f5ba83cfa80.md"
$file1Path  = "e2e\a41100b7-a95c-4c82-9f5a-4f5ba83cfa80.md"
$file1ZhXlf = "a41100b7-a95c-4c82-9f5a-4f5ba83cfa80.8ea73ab795626e51148c53dba48ba76366bc4054.zh-cn.xlf"
$file1DeXlf = "a41100b7-a95c-4c82-9f5a-4f5ba83cfa80.8ea73ab795626e51148c53dba48ba76366bc4054.de-de.xlf"
$file1Url   = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f8c82938a2e1f926a623c832059cd149bae8af50/e2e/a41100b7-a95c-4c82-9f5a-4f5ba83cfa80.md"

$file2Name  = "be004eef-2a37-4684-98e5-e3457d0576fd.md"
$file2Path  = "e2e\be004eef-2a37-4684-98e5-e3457d0576fd.md"
$file2ZhXlf = "be004eef-2a37-4684-98e5-e3457d0576fd.c110ae1e3cf41afb30949d776138891f57cde364.zh-cn.xlf"
$file2DeXlf = "be004eef-2a37-4684-98e5-e3457d0576fd.c110ae1e3cf41afb30949d776138891f57cde364.de-de.xlf"
$file2Url   = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f8c82938a2e1f926a623c832059cd149bae8af50/e2e/be004eef-2a37-4684-98e5-e3457d0576fd.md"

# ---------------------------------------------------------------------
# Sheet "Overview": columns A-G
#   A File Name | B Path And Name | C Extension | D Publish URL
#   E zh-cn | F de-de | G Latest HO Xliff Generate Date
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)

$rowA = $loOverview.ListRows.Add()
$wsOverview.Range("A6").Value = $file1Name
$wsOverview.Range("C6").Value = ".md"
$wsOverview.Range("D6").Value = ""
$wsOverview.Range("E6").Value = "Ready for handoff"
$wsOverview.Range("F6").Value = "Ready for handoff"
$wsOverview.Range("G6").Value = $dateHandoff
$wsOverview.Range("G6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B6"), $file1Url, "", "", $file1Path)

$rowB = $loOverview.ListRows.Add()
$wsOverview.Range("A7").Value = $file2Name
$wsOverview.Range("C7").Value = ".md"
$wsOverview.Range("D7").Value = ""
$wsOverview.Range("E7").Value = "Ready for handoff"
$wsOverview.Range("F7").Value = "Ready for handoff"
$wsOverview.Range("G7").Value = $dateHandoff
$wsOverview.Range("G7").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B7"), $file2Url, "", "", $file2Path)

# ---------------------------------------------------------------------
# Sheet "zh-cn": columns A-P
#   A Source File Name | B File Extension | C Status | D Source Path
#   E Priority | F Content Duplicate | G Latest Handoff File
#   H Latest Handoff Datetime | I Latest Target File | J Latest Handback File
#   K Latest Handback DateTime | L Reference Tokens | M To be localized
#   N Dependency From | O Has metadata | P Error Detail
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)

$rowC = $loZhCn.ListRows.Add()
$wsZhCn.Range("B6").Value = ".md"
$wsZhCn.Range("C6").Value = "Ready for handoff"
$wsZhCn.Range("D6").Value = "e2e"
$wsZhCn.Range("E6").Value = "ht"
$wsZhCn.Range("F6").Value = "False"
$wsZhCn.Range("G6").Value = $file1ZhXlf
$wsZhCn.Range("H6").Value = $dateZhCn
$wsZhCn.Range("H6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("I6").Value = ""
$wsZhCn.Range("J6").Value = ""
$wsZhCn.Range("K6").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("K6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("L6").Value = ""
$wsZhCn.Range("M6").Value = "True"
$wsZhCn.Range("N6").Value = ""
$wsZhCn.Range("O6").Value = "False"
$wsZhCn.Range("P6").Value = ""
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A6"), $file1Url, "", "", $file1Name)

$rowD = $loZhCn.ListRows.Add()
$wsZhCn.Range("B7").Value = ".md"
$wsZhCn.Range("C7").Value = "Ready for handoff"
$wsZhCn.Range("D7").Value = "e2e"
$wsZhCn.Range("E7").Value = "ht"
$wsZhCn.Range("F7").Value = "False"
$wsZhCn.Range("G7").Value = $file2ZhXlf
$wsZhCn.Range("H7").Value = $dateZhCn
$wsZhCn.Range("H7").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("I7").Value = ""
$wsZhCn.Range("J7").Value = ""
$wsZhCn.Range("K7").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("K7").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("L7").Value = ""
$wsZhCn.Range("M7").Value = "True"
$wsZhCn.Range("N7").Value = ""
$wsZhCn.Range("O7").Value = "False"
$wsZhCn.Range("P7").Value = ""
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A7"), $file2Url, "", "", $file2Name)

# ---------------------------------------------------------------------
# Sheet "de-de": same column layout as "zh-cn"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)

$rowE = $loDeDe.ListRows.Add()
$wsDeDe.Range("B6").Value = ".md"
$wsDeDe.Range("C6").Value = "Ready for handoff"
$wsDeDe.Range("D6").Value = "e2e"
$wsDeDe.Range("E6").Value = "ht"
$wsDeDe.Range("F6").Value = "False"
$wsDeDe.Range("G6").Value = $file1DeXlf
$wsDeDe.Range("H6").Value = $dateHandoff
$wsDeDe.Range("H6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("I6").Value = ""
$wsDeDe.Range("J6").Value = ""
$wsDeDe.Range("K6").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("K6").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("L6").Value = ""
$wsDeDe.Range("M6").Value = "True"
$wsDeDe.Range("N6").Value = ""
$wsDeDe.Range("O6").Value = "False"
$wsDeDe.Range("P6").Value = ""
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A6"), $file1Url, "", "", $file1Name)

$rowF = $loDeDe.ListRows.Add()
$wsDeDe.Range("B7").Value = ".md"
$wsDeDe.Range("C7").Value = "Ready for handoff"
$wsDeDe.Range("D7").Value = "e2e"
$wsDeDe.Range("E7").Value = "ht"
$wsDeDe.Range("F7").Value = "False"
$wsDeDe.Range("G7").Value = $file2DeXlf
$wsDeDe.Range("H7").Value = $dateHandoff
$wsDeDe.Range("H7").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("I7").Value = ""
$wsDeDe.Range("J7").Value = ""
$wsDeDe.Range("K7").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("K7").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("L7").Value = ""
$wsDeDe.Range("M7").Value = "True"
$wsDeDe.Range("N7").Value = ""
$wsDeDe.Range("O7").Value = "False"
$wsDeDe.Range("P7").Value = ""
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A7"), $file2Url, "", "", $file2Name)

Write-Host "Report rows added for handoff."
